$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: "Fraction of archaea"
#   Value column (B3): text "0.06" -> numeric 0.06
#   Uncertainty column (D3): text "18.2" -> numeric 15.8
$ws.Range("B3").Value = 0.06
$ws.Range("D3").Value = 15.8

# Row 4: "Fraction of bacteria"
#   Value column (B4): text "0.94" -> numeric 0.94
#   Uncertainty column (D4): text "1.4" -> numeric 1.6
$ws.Range("B4").Value = 0.94
$ws.Range("D4").Value = 1.6

# Row 2: "Total biomass of bacteria and archaea ..."
#   Uncertainty column (D2): numeric 10 -> text "20.0"
# Force text storage (otherwise the numeric-looking string would be
# re-interpreted as a number), then drop back to the default style so no
# lingering number-format is left applied to the cell.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "20.0"
$ws.Range("D2").Style = "Normal"
